$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 (a new Lasso Regression+normalization+lag1+PCA(2) result),
# shifting the existing rows 7-11 (Polynomial Regression ... GB Regression) down to 8-12.
$ws.Rows("7:7").Insert()

# Copy the formatting from the (now) row below into the freshly inserted row so it
# picks up the same borders/style as the other data rows.
$ws.Range("A8:C8").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row with the new model result.
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Lasso Regression+normalization+ lag1+PCA(2)"
$ws.Range("C7").Value = 81.644887362295094

# Renumber the Id column for the rows that were pushed down.
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10

# Widen column B slightly to fit the longer model names.
$ws.Columns("B:B").ColumnWidth = 43.28

# Update the selected cell to match the author's final selection.
$ws.Range("C6").Select()
